# Update "Waste Gas" capacity values (column E) for Open year 2022 and 2024
# to reflect updated data from upstream processes through 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24 -> Open year 2022: Waste Gas 206.94 -> 214.48
$ws.Range("E24").Value = 214.48

# Row 26 -> Open year 2024: Waste Gas 298.85 -> 443.665
$ws.Range("E26").Value = 443.665
